$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filterd")

# Update BASE AMOUNT (F), INITIAL AMOUNT (G), TOTAL (H) for rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = 2000

    if ($r -le 4) {
        $ws.Cells.Item($r, 7).Value = 5000
        $ws.Cells.Item($r, 8).Value = 7000
    } else {
        $ws.Cells.Item($r, 8).Value = 2000
    }
}
